# Applies the account.xlsx change:
#  - Adds a new "search_full_text" search-field row to the Search sheet
#    (right after "search_criteria" row, i.e. as the new row 4).
#  - Replaces the "account_homeAddress"/"homeAddress" row and the
#    "Role"/"roles" row with a single "Security Roles"/"securityRoles" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# 1) Insert a new row before current row 4 ("account_username" row),
#    shifting everything down by one.
$ws.Rows.Item(4).Insert()

# 2) Populate the newly inserted row 4 with the search_full_text entry.
$ws.Cells.Item(4, 1).Value = "`${msg.getProperty('search_full_text')}"
$ws.Cells.Item(4, 2).Value = "`${search_full_text}"

# After the insert, the old row 13 (account_homeAddress / homeAddress) is
# now row 14, and the old row 14 (Role / roles) is now row 15.

# 3) Delete the (now) row 14 that holds account_homeAddress / homeAddress,
#    which shifts the Role / roles row up from 15 to 14.
$ws.Rows.Item(14).Delete()

# 4) Overwrite the (now) row 14 -- previously Role / roles -- with the new
#    Security Roles / securityRoles entry.
$ws.Cells.Item(14, 1).Value = "`${msg.getProperty('Security Roles')}"
$ws.Cells.Item(14, 2).Value = "`${securityRoles}"
